# Update table title and column headers from October to November
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Table 6.2.C. Net Summer Capacity of Utility Scale Units Using Primarily Fossil Fuels and by State, November 2016 and 2015 (Megawatts)"

# "October 2016" -> "November 2016" header cells.
# NumberFormat is temporarily switched to text ("@") while assigning the
# value so Excel doesn't auto-convert the month/year text into a date
# serial number; the original number format is then restored so the
# cell's style stays identical to before.
foreach ($addr in @("B3","D3","F3","H3","J3","L3","N3","P3")) {
    $c = $ws.Range($addr)
    $fmt = $c.NumberFormat
    $c.NumberFormat = "@"
    $c.Value = "November 2016"
    $c.NumberFormat = $fmt
}

# "October 2015" -> "November 2015" header cells
foreach ($addr in @("C3","E3","G3","I3","K3","M3","O3","Q3")) {
    $c = $ws.Range($addr)
    $fmt = $c.NumberFormat
    $c.NumberFormat = "@"
    $c.Value = "November 2015"
    $c.NumberFormat = $fmt
}

# Updated data values (EPM_2016_11 run)
$ws.Range("B4").Value = 11911
$ws.Range("E4").Value = 1111.7
$ws.Range("G4").Value = 644.29999999999995
$ws.Range("P4").Value = 22730.1
$ws.Range("Q4").Value = 22753.3
$ws.Range("G5").Value = 419.1
$ws.Range("Q5").Value = 6309.6
$ws.Range("E7").Value = 331.5
$ws.Range("Q7").Value = 9839.2999999999993
$ws.Range("B9").Value = 1779.7
$ws.Range("P9").Value = 1809.3
$ws.Range("B11").Value = 26296.799999999999
$ws.Range("C11").Value = 24621.1
$ws.Range("E11").Value = 7599.7
$ws.Range("P11").Value = 69323.899999999994
$ws.Range("Q11").Value = 68015.600000000006
$ws.Range("C12").Value = 8047.5
$ws.Range("E12").Value = 2817.1
$ws.Range("Q12").Value = 13535.5
$ws.Range("B13").Value = 8122
$ws.Range("C13").Value = 8069.1
$ws.Range("P13").Value = 25983
$ws.Range("Q13").Value = 26310.5
$ws.Range("H15").Value = 62758.3
$ws.Range("P15").Value = 113805.9
$ws.Range("H20").Value = 7361.4
$ws.Range("P20").Value = 14340.1
$ws.Range("G21").Value = 3617.1
$ws.Range("I21").Value = 36198.1
$ws.Range("M21").Value = 4096
$ws.Range("Q21").Value = 61263.3
$ws.Range("G23").Value = 2024
$ws.Range("I23").Value = 4687.2
$ws.Range("M23").Value = 538.1
$ws.Range("Q23").Value = 9583.1
$ws.Range("M24").Value = 799.4
$ws.Range("Q24").Value = 10240.6
$ws.Range("C29").Value = 48550.2
$ws.Range("F29").Value = 7112.1
$ws.Range("G29").Value = 7157.8
$ws.Range("H29").Value = 58281.599999999999
$ws.Range("I29").Value = 59261.599999999999
$ws.Range("M29").Value = 11981.4
$ws.Range("Q29").Value = 158279.9
$ws.Range("C32").Value = 26245.599999999999
$ws.Range("G32").Value = 3119.7
$ws.Range("M32").Value = 5936.9
$ws.Range("Q32").Value = 53374.6
$ws.Range("M36").Value = 525.4
$ws.Range("Q36").Value = 11635.2
$ws.Range("I37").Value = 4609.3
$ws.Range("Q37").Value = 16845.8
$ws.Range("F38").Value = 123
$ws.Range("H38").Value = 12958
$ws.Range("I38").Value = 13081
$ws.Range("Q38").Value = 14163.3
$ws.Range("I39").Value = 31803.4
$ws.Range("Q39").Value = 68051.3
$ws.Range("I41").Value = 13436.7
$ws.Range("Q41").Value = 19153.5
$ws.Range("D58").Value = 11949.7
$ws.Range("E58").Value = 11509.2
$ws.Range("G58").Value = 12668.7
$ws.Range("P58").Value = 52496.4
$ws.Range("Q58").Value = 51992.9
$ws.Range("D59").Value = 11054.5
$ws.Range("E59").Value = 10614
$ws.Range("G59").Value = 12416.7
$ws.Range("P59").Value = 43308.9
$ws.Range("Q59").Value = 43325.4
$ws.Range("L62").Value = 2591.5
$ws.Range("P62").Value = 4182.2
$ws.Range("L63").Value = 723.9
$ws.Range("P63").Value = 2125
$ws.Range("B65").Value = 239449.60000000001
$ws.Range("C65").Value = 234018
$ws.Range("D65").Value = 125081.2
$ws.Range("E65").Value = 123605.6
$ws.Range("F65").Value = 83428.3
$ws.Range("G65").Value = 82844.399999999994
$ws.Range("H65").Value = 271591.90000000002
$ws.Range("I65").Value = 281254.2
$ws.Range("L65").Value = 34875.599999999999
$ws.Range("M65").Value = 35101
$ws.Range("P65").Value = 758464.1
$ws.Range("Q65").Value = 761097.9
